$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = 112043819
$ws.Range("B2").Value = 78699
$ws.Range("E2").Value = 6458
$ws.Range("F2").Value = "Lunglav"
$ws.Range("G2").Value = "Lobaria pulmonaria"
$ws.Range("H2").Value = "(L.) Hoffm."
$ws.Range("K2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("Q2").Value = 547979
$ws.Range("R2").Value = 6960195
$ws.Range("AC2").ClearContents()

# --- Row 3 updates ---
$ws.Range("A3").Value = 112043839
$ws.Range("B3").Value = 78699
$ws.Range("Q3").Value = 547969
$ws.Range("R3").Value = 6960405

# --- Row 4 updates ---
$ws.Range("A4").Value = 112043807
$ws.Range("B4").Value = 56430
$ws.Range("E4").Value = 100109
$ws.Range("F4").Value = "Tretåig hackspett"
$ws.Range("G4").Value = "Picoides tridactylus"
$ws.Range("H4").Value = "(Linnaeus, 1758)"
$ws.Range("Q4").Value = 547961
$ws.Range("R4").Value = 6960421
$ws.Range("AC4").Value = "ringhack"
